$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 21031.2
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 21031.2
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H98").Value = 2185.5
$ws.Range("I98").Value = 2247.8948
$ws.Range("K98").Value = 2247.8948
$ws.Range("M98").Value = -749.8948
$ws.Range("H121").Value = 998
$ws.Range("J121").Value = 998
$ws.Range("L121").Value = 2994
$ws.Range("N121").Value = -6488
$ws.Range("H122").Value = 2185.5
$ws.Range("I122").Value = 2247.8948
$ws.Range("K122").Value = 6743.6844
$ws.Range("M122").Value = -4293.6844
$ws.Range("H129").Value = 918.34784
$ws.Range("J129").Value = 892.6512
$ws.Range("L129").Value = 2677.9536
$ws.Range("N129").Value = -12677.9536
$ws.Range("H132").Value = 1106.5588
$ws.Range("I132").Value = 1069.7587
$ws.Range("J132").Value = 1320
$ws.Range("K132").Value = 3209.2761
$ws.Range("L132").Value = 3960
$ws.Range("M132").Value = -679.2761
$ws.Range("N132").Value = -9020
$ws.Range("H137").Value = 1378.091
$ws.Range("I137").Value = 1282.375
$ws.Range("J137").Value = 1633.3334
$ws.Range("K137").Value = 3847.125
$ws.Range("L137").Value = 4900.0002
$ws.Range("M137").Value = -1297.125
$ws.Range("N137").Value = -10000.0002
$ws.Range("H138").Value = 2871.17
$ws.Range("J138").Value = 2811.7144
$ws.Range("L138").Value = 8435.143199999999
$ws.Range("N138").Value = -18715.1432
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3948.6938
$ws.Range("I32").Value = 2797.561
$ws.Range("J32").Value = 9848.25
$ws.Range("K32").Value = 2797.561
$ws.Range("L32").Value = 9848.25
$ws.Range("M32").Value = -2510.561
$ws.Range("N32").Value = -10422.25
$ws.Range("H61").Value = 4627.533
$ws.Range("I61").Value = 3227
$ws.Range("J61").Value = 10229.667
$ws.Range("K61").Value = 3227
$ws.Range("L61").Value = 10229.667
$ws.Range("M61").Value = -3015
$ws.Range("N61").Value = -10653.667
$ws.Range("H63").Value = 6200.8
$ws.Range("I63").Value = 6200.8
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 6200.8
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -5514.8
$ws.Range("H66").Value = 6200.8
$ws.Range("I66").Value = 6200.8
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 31004
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -27572
$ws.Range("H74").Value = 1098.8
$ws.Range("J74").Value = 2705.111
$ws.Range("L74").Value = 2705.111
$ws.Range("N74").Value = -4453.111
$ws.Range("H77").Value = 1098.8
$ws.Range("J77").Value = 2705.111
$ws.Range("L77").Value = 13525.555
$ws.Range("N77").Value = -22261.555
$ws.Range("H122").Value = 2012
$ws.Range("I122").Value = 2012
$ws.Range("K122").Value = 6036
$ws.Range("M122").Value = -3586
$ws.Range("H132").Value = 1576.3793
$ws.Range("I132").Value = 968.7368
$ws.Range("K132").Value = 2906.2104
$ws.Range("M132").Value = -376.2103999999999
$ws.Range("H136").Value = 4627.533
$ws.Range("I136").Value = 3227
$ws.Range("J136").Value = 10229.667
$ws.Range("K136").Value = 9681
$ws.Range("L136").Value = 30689.001
$ws.Range("M136").Value = -7131
$ws.Range("N136").Value = -35789.001
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1669.8
$ws.Range("I99").Value = 1385.4286
$ws.Range("K99").Value = 1385.4286
$ws.Range("M99").Value = 112.5714
$ws.Range("H134").Value = 15743.611
$ws.Range("I134").Value = 15956.143
$ws.Range("K134").Value = 47868.429
$ws.Range("M134").Value = -45333.429

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2711.5366
$ws.Range("J31").Value = 3974.5625
$ws.Range("L31").Value = 3974.5625
$ws.Range("N31").Value = -4564.5625
$ws.Range("H34").Value = 2711.5366
$ws.Range("J34").Value = 3974.5625
$ws.Range("L34").Value = 3974.5625
$ws.Range("N34").Value = -4378.5625
$ws.Range("H58").Value = 2416892.5
$ws.Range("I58").Value = 3345589
$ws.Range("K58").Value = 3345589
$ws.Range("M58").Value = -3345386
$ws.Range("H132").Value = 2003.0358
$ws.Range("I132").Value = 1113.7
$ws.Range("K132").Value = 3341.1
$ws.Range("M132").Value = -811.1000000000004
$ws.Range("H134").Value = 741.76666
$ws.Range("I134").Value = 741.76666
$ws.Range("K134").Value = 2225.29998
$ws.Range("M134").Value = 309.7000200000002
$ws.Range("H136").Value = 2416892.5
$ws.Range("I136").Value = 3345589
$ws.Range("K136").Value = 10036767
$ws.Range("M136").Value = -10034217

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 119.04348
$ws.Range("I4").Value = 119.04348
$ws.Range("K4").Value = 357.13044
$ws.Range("M4").Value = -245.13044
$ws.Range("H17").Value = 333346660
$ws.Range("I17").Value = 1000000000
$ws.Range("J17").Value = 20001
$ws.Range("K17").Value = 3000000000
$ws.Range("L17").Value = 60003
$ws.Range("M17").Value = -2999999831
$ws.Range("N17").Value = -60341
$ws.Range("H33").Value = 130.14285
$ws.Range("J33").Value = 164.5
$ws.Range("L33").Value = 987
$ws.Range("N33").Value = -1553
$ws.Range("H88").Value = 4999.5
$ws.Range("J88").Value = 5399.4
$ws.Range("L88").Value = 16198.2
$ws.Range("N88").Value = -17054.2
$ws.Range("H91").Value = 4999.5
$ws.Range("J91").Value = 5399.4
$ws.Range("L91").Value = 16198.2
$ws.Range("N91").Value = -19162.2
$ws.Range("H131").Value = 9013.319
$ws.Range("I131").Value = 564
$ws.Range("J131").Value = 9877.454
$ws.Range("K131").Value = 1692
$ws.Range("L131").Value = 29632.362
$ws.Range("M131").Value = 3348
$ws.Range("N131").Value = -39712.362

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 63031
$ws.Range("J47").Value = 63031
$ws.Range("L47").Value = 63031
$ws.Range("N47").Value = -64167
$ws.Range("H97").Value = 1262.6666
$ws.Range("I97").Value = 1149
$ws.Range("J97").Value = 1575.25
$ws.Range("K97").Value = 1149
$ws.Range("L97").Value = 1575.25
$ws.Range("M97").Value = -653
$ws.Range("N97").Value = -2567.25
$ws.Range("H132").Value = 2028389.2
$ws.Range("I132").Value = 2567160
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 7701480
$ws.Range("L132").Value = 23997
$ws.Range("M132").Value = -7698950
$ws.Range("N132").Value = -29057
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9428.286
$ws.Range("I40").Value = 1999.6666
$ws.Range("K40").Value = 1999.6666
$ws.Range("M40").Value = -1863.6666
$ws.Range("H100").Value = 1130.125
$ws.Range("I100").Value = 1092.2
$ws.Range("J100").Value = 1193.3334
$ws.Range("K100").Value = 1092.2
$ws.Range("L100").Value = 1193.3334
$ws.Range("M100").Value = -551.2
$ws.Range("N100").Value = -2275.3334
$ws.Range("H132").Value = 3036.6
$ws.Range("I132").Value = 2399.7144
$ws.Range("K132").Value = 7199.1432
$ws.Range("M132").Value = -4669.1432
$ws.Range("H136").Value = 3243.6365
$ws.Range("I136").Value = 1487.5454
$ws.Range("K136").Value = 4462.6362
$ws.Range("M136").Value = -1912.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3213.9443
$ws.Range("I132").Value = 2715.818
$ws.Range("J132").Value = 3996.7144
$ws.Range("K132").Value = 8147.454000000001
$ws.Range("L132").Value = 11990.1432
$ws.Range("M132").Value = -5617.454000000001
$ws.Range("N132").Value = -17050.1432
$ws.Range("H136").Value = 17363560
$ws.Range("I136").Value = 29242164
$ws.Range("K136").Value = 87726492
$ws.Range("M136").Value = -87723942
